# Re-creates the author's edit:
#   1. The table on slide 16 gets a new table style
#      ({1505EEF1-7193-4265-A812-C80C956F96D1} -> {5D916D01-7059-4B79-BD5C-1ACDF8F195F8}).
#   2. The deck's theme colour scheme (physically stored in ppt/theme/theme2.xml,
#      the theme actually wired to the slide master) is switched from the
#      "Integral" palette to the standard "Office" palette.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------
$tableSlide = $p.Slides.Item(16)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{5D916D01-7059-4B79-BD5C-1ACDF8F195F8}")
    }
}

# --- 2. Theme colours: Integral -> Office ----------------------------------
# Order exposed by ThemeColorScheme: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeColors = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $hex = $officeColors[$i - 1]
    $r = ($hex -shr 16) -band 0xFF
    $g = ($hex -shr 8) -band 0xFF
    $b = $hex -band 0xFF
    $ole = $r + ($g * 256) + ($b * 65536)
    $themeColors.Item($i).RGB = $ole
}
